$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price and volume(1h) figures per the latest data refresh.
# D-column price values are prefixed with a literal apostrophe so Excel keeps them
# as text (matching the inlineStr cells in the workbook) instead of auto-converting
# numeric-looking strings (e.g. "1.00", "7.60") into numbers.
$ws.Range("D2").Value = "68.331.15"
$ws.Range("E2").Value = "  -0.11%  "
$ws.Range("D3").Value = "2.707.24"
$ws.Range("E3").Value = "  +2.24%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'607.61"
$ws.Range("E5").Value = "  +1.85%  "
$ws.Range("D6").Value = "'166.63"
$ws.Range("E6").Value = "  +4.84%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "'0.559"
$ws.Range("E8").Value = "  +3.54%  "
$ws.Range("D9").Value = "2.707.93"
$ws.Range("E9").Value = "  +2.31%  "
$ws.Range("D10").Value = "'0.145"
$ws.Range("E10").Value = "  +1.73%  "
$ws.Range("E11").Value = "  +0.60%  "
$ws.Range("E12").Value = "  +3.19%  "
$ws.Range("D13").Value = "'5.29"
$ws.Range("E13").Value = "  +0.26%  "
$ws.Range("D14").Value = "'28.35"
$ws.Range("E14").Value = "  +1.21%  "
$ws.Range("D15").Value = "3.202.24"
$ws.Range("E15").Value = "  +2.26%  "
$ws.Range("E16").Value = "  -0.42%  "
$ws.Range("D17").Value = "68.266.74"
$ws.Range("E17").Value = "  -0.08%  "
$ws.Range("D18").Value = "2.707.80"
$ws.Range("E18").Value = "  +2.68%  "
$ws.Range("D19").Value = "'11.78"
$ws.Range("E19").Value = "  +1.99%  "
$ws.Range("D20").Value = "'369.69"
$ws.Range("E20").Value = "  +1.73%  "
$ws.Range("D21").Value = "'7.60"
$ws.Range("E21").Value = "  +1.81%  "
$ws.Range("D22").Value = "'4.47"
$ws.Range("E22").Value = "  +1.51%  "
$ws.Range("D23").Value = "'4.93"
$ws.Range("E23").Value = "  +3.27%  "
$ws.Range("E24").Value = "  -0.51%  "
$ws.Range("D25").Value = "'73.09"
$ws.Range("E25").Value = "  -2.18%  "
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("D27").Value = "'9.97"
$ws.Range("E27").Value = "  +0.83%  "
$ws.Range("E29").Value = "  +0.30%  "
$ws.Range("E30").Value = "  -0.18%  "
$ws.Range("D31").Value = "'577.92"
$ws.Range("E31").Value = "  +2.12%  "
$ws.Range("D32").Value = "'8.09"
$ws.Range("E32").Value = "  +0.31%  "
$ws.Range("D33").Value = "'1.41"
$ws.Range("E33").Value = "  +0.91%  "
$ws.Range("E34").Value = "  +5.37%  "
$ws.Range("E35").Value = "  +1.73%  "
$ws.Range("D36").Value = "'1.00"
$ws.Range("E36").Value = "  +0.04%  "
$ws.Range("E37").Value = "  -3.47%  "
$ws.Range("D38").Value = "'161.36"
$ws.Range("E38").Value = "  +0.42%  "
$ws.Range("D39").Value = "'19.81"
$ws.Range("E39").Value = "  +0.80%  "
$ws.Range("D40").Value = "'0.376"
$ws.Range("E40").Value = "  +1.66%  "
$ws.Range("D41").Value = "'1.86"
$ws.Range("E41").Value = "  -0.32%  "
$ws.Range("D42").Value = "'5.36"
$ws.Range("E42").Value = "  +0.72%  "
$ws.Range("D43").Value = "'17.98"
$ws.Range("E43").Value = "  +0.96%  "
$ws.Range("E44").Value = "  -1.87%  "
$ws.Range("E45").Value = "  -0.05%  "
$ws.Range("E46").Value = "  -3.38%  "
$ws.Range("D47").Value = "'40.76"
$ws.Range("E47").Value = "  +1.09%  "
$ws.Range("D48").Value = "'0.594"
$ws.Range("E48").Value = "  +3.40%  "
$ws.Range("D49").Value = "'154.24"
$ws.Range("E49").Value = "  -2.53%  "
$ws.Range("E50").Value = "  +1.72%  "
$ws.Range("D51").Value = "'1.76"
$ws.Range("E51").Value = "  +3.87%  "
